$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1339.2858
$ws.Range("I18").Value = 1225
$ws.Range("K18").Value = 1225
$ws.Range("M18").Value = -941

$ws.Range("H33").Value = 382.03226
$ws.Range("I33").Value = 354.07693
$ws.Range("K33").Value = 354.07693
$ws.Range("M33").Value = -125.07693

$ws.Range("H70").Value = 1744.5
$ws.Range("I70").Value = 1753.3334
$ws.Range("J70").Value = 1740.7142
$ws.Range("K70").Value = 5260.0002
$ws.Range("L70").Value = 5222.142599999999
$ws.Range("M70").Value = -4990.0002
$ws.Range("N70").Value = -5762.142599999999

$ws.Range("H73").Value = 1744.5
$ws.Range("I73").Value = 1753.3334
$ws.Range("J73").Value = 1740.7142
$ws.Range("K73").Value = 5260.0002
$ws.Range("L73").Value = 5222.142599999999
$ws.Range("M73").Value = -4324.0002
$ws.Range("N73").Value = -7094.142599999999

$ws.Range("H133").Value = 35594
$ws.Range("J133").Value = 35594
$ws.Range("L133").Value = 35594
$ws.Range("N133").Value = -45714

$ws.Range("H137").Value = 1535.7576
$ws.Range("I137").Value = 1053.4706
$ws.Range("J137").Value = 2048.1875
$ws.Range("K137").Value = 3160.4118
$ws.Range("L137").Value = 6144.5625
$ws.Range("M137").Value = -610.4118000000003
$ws.Range("N137").Value = -11244.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6968.1035
$ws.Range("I32").Value = 6968.1035
$ws.Range("K32").Value = 6968.1035
$ws.Range("M32").Value = -6681.1035

$ws.Range("H45").Value = 1159.45
$ws.Range("I45").Value = 1050.7333
$ws.Range("J45").Value = 1485.6
$ws.Range("K45").Value = 1050.7333
$ws.Range("L45").Value = 1485.6
$ws.Range("M45").Value = -673.7333000000001
$ws.Range("N45").Value = -2239.6

$ws.Range("H74").Value = 1389.8
$ws.Range("I74").Value = 730.625
$ws.Range("K74").Value = 730.625
$ws.Range("M74").Value = 143.375

$ws.Range("H77").Value = 1389.8
$ws.Range("I77").Value = 730.625
$ws.Range("K77").Value = 3653.125
$ws.Range("M77").Value = 714.875

$ws.Range("H101").Value = 34666.332
$ws.Range("J101").Value = 34666.332
$ws.Range("L101").Value = 34666.332
$ws.Range("N101").Value = -41156.332

$ws.Range("H122").Value = 1452.5714
$ws.Range("I122").Value = 1006
$ws.Range("J122").Value = 2345.7144
$ws.Range("K122").Value = 3018
$ws.Range("L122").Value = 7037.1432
$ws.Range("M122").Value = -568
$ws.Range("N122").Value = -11937.1432

$ws.Range("H132").Value = 3025.037
$ws.Range("I132").Value = 2607.1428
$ws.Range("K132").Value = 7821.428400000001
$ws.Range("M132").Value = -5291.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 50001080
$ws.Range("I99").Value = 52632624
$ws.Range("K99").Value = 52632624
$ws.Range("M99").Value = -52631126

$ws.Range("H102").Value = 14040.917
$ws.Range("I102").Value = 5943.4443
$ws.Range("K102").Value = 5943.4443
$ws.Range("M102").Value = -2698.4443

$ws.Range("H134").Value = 3297.682
$ws.Range("I134").Value = 750.30304
$ws.Range("J134").Value = 10939.818
$ws.Range("K134").Value = 2250.90912
$ws.Range("L134").Value = 32819.454
$ws.Range("M134").Value = 284.0908799999997
$ws.Range("N134").Value = -37889.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1593.3636
$ws.Range("I31").Value = 1837.0769
$ws.Range("J31").Value = 1434.95
$ws.Range("K31").Value = 1837.0769
$ws.Range("L31").Value = 1434.95
$ws.Range("M31").Value = -1542.0769
$ws.Range("N31").Value = -2024.95

$ws.Range("H34").Value = 1593.3636
$ws.Range("I34").Value = 1837.0769
$ws.Range("J34").Value = 1434.95
$ws.Range("K34").Value = 1837.0769
$ws.Range("L34").Value = 1434.95
$ws.Range("M34").Value = -1635.0769
$ws.Range("N34").Value = -1838.95

$ws.Range("H86").Value = 6100108
$ws.Range("I86").Value = 16694416
$ws.Range("J86").Value = 46217.715
$ws.Range("K86").Value = 16694416
$ws.Range("L86").Value = 46217.715
$ws.Range("M86").Value = -16693293
$ws.Range("N86").Value = -48463.715

$ws.Range("H89").Value = 6100108
$ws.Range("I89").Value = 16694416
$ws.Range("J89").Value = 46217.715
$ws.Range("K89").Value = 83472080
$ws.Range("L89").Value = 231088.575
$ws.Range("M89").Value = -83466464
$ws.Range("N89").Value = -242320.575

$ws.Range("H132").Value = 2014.6842
$ws.Range("I132").Value = 1485.4
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 4456.200000000001
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -1926.200000000001
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 231.05556
$ws.Range("I14").Value = 231.05556
$ws.Range("K14").Value = 693.16668
$ws.Range("M14").Value = -520.16668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("N26").Value = -25560

$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -25996

$ws.Range("H70").Value = 34618124
$ws.Range("I70").Value = 31252452
$ws.Range("J70").Value = 40003200
$ws.Range("K70").Value = 31252452
$ws.Range("L70").Value = 40003200
$ws.Range("M70").Value = -31252182
$ws.Range("N70").Value = -40003740

$ws.Range("H73").Value = 34618124
$ws.Range("I73").Value = 31252452
$ws.Range("J73").Value = 40003200
$ws.Range("K73").Value = 31252452
$ws.Range("L73").Value = 40003200
$ws.Range("M73").Value = -31251516
$ws.Range("N73").Value = -40005072

$ws.Range("H80").Value = 5050.0835
$ws.Range("I80").Value = 4625
$ws.Range("J80").Value = 5900.25
$ws.Range("K80").Value = 4625
$ws.Range("L80").Value = 5900.25
$ws.Range("M80").Value = -3627
$ws.Range("N80").Value = -7896.25

$ws.Range("H83").Value = 5050.0835
$ws.Range("I83").Value = 4625
$ws.Range("J83").Value = 5900.25
$ws.Range("K83").Value = 23125
$ws.Range("L83").Value = 29501.25
$ws.Range("M83").Value = -18133
$ws.Range("N83").Value = -39485.25

$ws.Range("H132").Value = 2523.4375
$ws.Range("I132").Value = 2483.7
$ws.Range("J132").Value = 2589.6667
$ws.Range("K132").Value = 7451.099999999999
$ws.Range("L132").Value = 7769.000100000001
$ws.Range("M132").Value = -4921.099999999999
$ws.Range("N132").Value = -12829.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1099.2
$ws.Range("I61").Value = 1124.125
$ws.Range("K61").Value = 1124.125
$ws.Range("M61").Value = -922.125

$ws.Range("H82").Value = 2112.7144
$ws.Range("I82").Value = 1964.8334
$ws.Range("K82").Value = 1964.8334
$ws.Range("M82").Value = -1603.8334

$ws.Range("H85").Value = 2112.7144
$ws.Range("I85").Value = 1964.8334
$ws.Range("K85").Value = 1964.8334
$ws.Range("M85").Value = -716.8334

$ws.Range("H113").Value = 1099.2
$ws.Range("I113").Value = 1124.125
$ws.Range("K113").Value = 1124.125
$ws.Range("M113").Value = 1045.875

$ws.Range("H132").Value = 2604.111
$ws.Range("J132").Value = 2957.3333
$ws.Range("L132").Value = 8871.999899999999
$ws.Range("N132").Value = -13931.9999

$ws.Range("H136").Value = 1731.2106
$ws.Range("I136").Value = 1176.4615
$ws.Range("K136").Value = 3529.3845
$ws.Range("M136").Value = -979.3844999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 240
$ws.Range("I17").Value = 240
$ws.Range("K17").Value = 240
$ws.Range("M17").Value = -68

$ws.Range("H136").Value = 1028.4117
$ws.Range("I136").Value = 1007.32
$ws.Range("J136").Value = 1087
$ws.Range("K136").Value = 3021.96
$ws.Range("L136").Value = 3261
$ws.Range("M136").Value = -471.96
$ws.Range("N136").Value = -8361

Write-Host "Applied Kujata_Profits updates"
